$p = $ppt.ActivePresentation
$sm = $p.SlideMaster
Write-Host "SlideMaster:" $sm
$hm = $p.HandoutMaster
Write-Host "HandoutMaster:" $hm
$hmTheme = $hm.Theme
$cs = $hmTheme.ThemeColorScheme
Write-Host "Handout theme accent1:" $cs.Item(5).RGB
